$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in H1, copying the formatting (xfId style) used by
# the other header cells (e.g. G1: bold, bordered, centered) before setting
# its text so the paste doesn't clobber the value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" column values for each existing data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
